$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")
$ws.Activate()

# Clear the credential/URL test data that was pulled from the workbook
# (userid, password, and the URL value in Q2/R2/S2), while leaving the
# styled-but-empty Q2 cell (keeps its "Hyperlink" style s="2").
$ws.Range("Q2:S2").ClearContents()

# Drop the now-stale hyperlink that pointed at the removed URL.
$ws.Hyperlinks.Delete()

# Reflect the new selection left in the sheet after the edit.
$ws.Range("Q2:S2").Select()
